$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.060031414031982
$ws.Range("B1").Value = 2.34818696975708
$ws.Range("C1").Value = 2.373142719268799
$ws.Range("D1").Value = 2.840351819992065
$ws.Range("E1").Value = 0.9352903962135315
